$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Applies every cell text update from the commit diff. D-column values that parse as
# plain numbers ("1.002", "45.18", ...) are round-tripped through a temporary text
# NumberFormat so the COM layer stores them as inline/shared strings (matching the
# original inlineStr cells) instead of silently coercing them to numeric <v> cells;
# the per-cell Style reset back to "Normal" keeps the cell styling unchanged afterwards.

$ws.Range("D2").Value = "30.006.12"
$ws.Range("E2").Value = "  +9.69%  "
$ws.Range("D3").Value = "1.867.37"
$ws.Range("E3").Value = "  +6.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4957"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2829"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06517"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.72%  "
$ws.Range("D11").Value = "1.868.30"
$ws.Range("E11").Value = "  +6.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.39%  "
$ws.Range("E13").Value = "  +3.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6590"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "85.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +10.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.793"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.15%  "
$ws.Range("D17").Value = "29.975.67"
$ws.Range("E17").Value = "  +9.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007455"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.36%  "
$ws.Range("E20").Value = "  +10.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "2.110.13"
$ws.Range("E22").Value = "  +7.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.713"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.06%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.988"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.56%  "
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.496"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +24.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.932"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.10%  "
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.230"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08567"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.868"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05047"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.126"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6798"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.694"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.326"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.67%  "
$ws.Range("E39").Value = "  +6.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9568"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01628"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.122"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("E44").Value = "  +2.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4155"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.377"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1245"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05628"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.285"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3696"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.51%  "
